$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2: host_organization N/A -> Elsevier BV
$ws.Range("G2").Value = "Elsevier BV"

# A new publication record (OpenAlex W4283017621, "Sa1710: ...") was inserted
# as row 3, and the pre-existing "Contributors" record (OpenAlex W4206641953)
# that used to be row 3 is now row 4. Concretely this swaps every column that
# differs between what used to be row 3 and row 4; identical columns (I, J, K,
# M, N, O) are left as-is.

$ws.Range("A3").Value = "Radhika Babaria, Hemanth Gavini, Cynthia J Doane, Chengcheng Hu, Bhaskar Banerjee"
$ws.Range("A4").Value = "Krishnendu Adhikary, Amit Agarwal, Md. Akil Hossain, Mohd Fadhlizil Fasihi Mohd Aluwi, Lin Ai, Siddaraju Anusha, Jess Armine, Radhika Babaria, Debasis Bagchi, Manashi Bagchi, Bhaskar Banerjee, Pradipta Banerjee, Samudra Prosad Banik, Tejaswini Baral, P. Basu, Bharathi Bethapudi, Jhimli Bhatttacharyya, Nabendu Biswas, Małgorzata Bukowiecka‐Matusiak, Izabela Burzyńska-Pędziwiatr, Leah Bush, Sanjoy Chakraborty, Ankita Chatterjee, Aritra Chatterjee, Sabyasachi Chatterjee, Smriti Chawla, Amitava Das, Dolan Das, Sujit Das, S. Devaraja, Bernard W. Downs, Jaclyn Downs, Nandini Ghosh, Оксана Головинская, Osamu Handa, Subrota Hati, Annaelle Hip Kam, Komal Jalan, Pradeep Kathi, Myung‐Sunny Kim, Aneta Kopeć, Subrahmanya Kumar Kukkupuni, Shilia Jacob Kurian, Steve Kushner, Hye Won Lee, Myeong Soo Lee, Himangshu Sekhar Maji, Labonya Mandal, A. Mavani, A.K.M. Moyeenul Huq, Fatima Muili, Deepak Mundkinajeddu, Sasikumar Murugan, Sreejayan Nair, Yuji Naito, Vidushi S. Neergheen, Pradeep Singh Negi, Kalu Ngele, Muruganantham Nithyanantham, Ya Fatou Njie‐Mbye, Sunny E. Ohia, Anthonia Okolie, Catherine A. Opere, Harry G. Preuss, Mahadev Rao, Moumita Ray, Akanksha Rout, Marufa Rumman, Hephzibah Saji, Saptadip Samanta, Riya Sarkar, Kenji Satô, Shalini Sehgal, Sonal Sekhar Miraj, M.N. Sharath kumar, Saki Shirako, Abhilasha Singh, Vineet Singh, Derek Smith, Eunhye Song, G. Sowmyashree, Tomohisa Takagi, Md. Hafiz Uddin, Chethala N. Vishnuprasad, Satoshi Wada, Chin‐Kun Wang, Lucyna A. Woźniak, Orie Yoshinari, Jerzy Zawistowski"

$ws.Range("B3").Value = "; ; ; ; "
$ws.Range("B4").Value = "Department of Interdisciplinary Sciences, Centurion University of Technology and Management, R. Sitapur, Odisha, India; Research and Development Center, Natural Remedies Private Limited, Bengaluru, Karnataka, India; Department of Pharmacology and Experimental Therapeutics, School of Medicine, Boston University, Boston, MA, United States; Faculty of Industrial Sciences and Technology, Universiti Malaysia Pahang, Lebuhraya Tun Razak, Gambang, Pahang, Malaysia; Korea Institute of Oriental Medicine, Daejeon, Korea; Korea University of Science and Techonology, Daejeon, Korea; Academy of Scientific and Innovative Research, Ghaziabad, CSIR-Central Food Technological Research Institute, Mysuru, Karnataka, India; Nutrigenomics and Functional Medicine, The Center for Bioindividualized Medicine, Hyannis, MA, United States; Division of Gastroenterology, Department of Internal Medicine, University of Arizona, Tucson, AZ, United States; College of Pharmacy and Health Sciences, Texas Southern University, Houston, TX, United States; Department of Biology, Adelphi University, Garden City, NY, United States; Department of R&D, Victory Nutrition Inc., Bonita Springs, FL, United States; Dr. Herbs LLC, Concord, CA, United States; Division of Gastroenterology, Department of Internal Medicine, University of Arizona, Tucson, AZ, United States; Department of Biochemistry and Plant Physiology, Centurion University of Technology and Management, R. Sitapur, Odisha, India; Department of Microbiology, Maulana Azad College, Kolkata, West Bengal, India; Department of Pharmacy Practice, Manipal College of Pharmaceutical Sciences, Manipal Academy of Higher Education, Manipal, Karnataka, India; National Institute of Biomedical Genomics, Kalyani, West Bengal, India; Research and Development Center, Natural Remedies Private Limited, Bengaluru, Karnataka, India; Department of Chemistry, National Institute of Technology Nagaland, Dimapur, Nagaland, India; Department of Life Sciences, Presidency University, Kolkata, West Bengal, India; Medical University of Lodz, Department of Structural Biology, Lodz, Poland; Medical University of Lodz, Department of Structural Biology, Lodz, Poland; ; ; ; ; ; ; ; ; ; Medical University of Lodz, Department of Structural Biology, Lodz, Poland; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; ; National Institute of Biomedical Genomics, Kalyani, West Bengal, India; ; ; ; ; ; ; ; ; ; Medical University of Lodz, Department of Structural Biology, Lodz, Poland; ; ; ; ; ; ; ; "

$ws.Range("C3").Value = "https://openalex.org/W4283017621"
$ws.Range("C4").Value = "https://openalex.org/W4206641953"

$ws.Range("D3").Value = "Sa1710: A MINIATURIZED MULTI-VIEW IMAGING DEVICE (MVID) FOR SIMULTANEOUS FORWARD AND REAR VIEWS IN COLONOSCOPY; FIRST IN-VIVO USE."
$ws.Range("D4").Value = "Contributors"

$ws.Range("E3").Value = "'2022-05-01"
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'2022-01-01"
$ws.Range("E4").Style = "Normal"

$ws.Range("F3").Value = "Gastroenterology"
$ws.Range("F4").Value = "Elsevier eBooks"

$ws.Range("G3").Value = "Elsevier BV"
$ws.Range("G4").Value = "N/A"

$ws.Range("H3").Value = "https://doi.org/10.1016/s0016-5085(22)61122-3"
$ws.Range("H4").Value = "https://doi.org/10.1016/b978-0-12-821232-5.01002-8"

$ws.Range("L3").Value = "de"
$ws.Range("L4").Value = "N/A"

$ws.Range("P3").Value = "https://doi.org/10.1016/s0016-5085(22)61122-3"
$ws.Range("P4").Value = "https://doi.org/10.1016/b978-0-12-821232-5.01002-8"

$ws.Range("Q3").Value = "article"
$ws.Range("Q4").Value = "book-chapter"

